$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.286832544864788, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 4.23731228292506)
    3  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    4  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    5  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 15.88780690183548)
    6  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    7  = @(0.1190320826869504, 0.306821227259698, 261.3203778131603, 10.19245300693656, 271.9386841300435)
    8  = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 10.19245300693656, 16.84135478251809)
    9  = @(0.2917716402565462, 0.04071648406533734, 3.537761648806719, 10.19245300693656, 14.06270278006516)
    10 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    11 = @(3.286832544864788, 117.745847958593, 0.7527432677738641, 10.19245300693656, 131.9778767781682)
    12 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    13 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 15.88780690183548)
    14 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    15 = @(3.286832544864788, 3286.919754855326, 22.3905356188092, 10.19245300693656, 3322.789576025937)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
